$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The rows for "Fruta, Vega Monumental Concepción - Caqui" report an updated
# weekly entry: swap the date/price/origin/unit data between row 2 and row 3
# (the descriptive columns A-M are identical on both rows).

# Save current (pre-edit) values of row 2's changing columns
$d2 = $ws.Range("D2").Value2
$n2 = $ws.Range("N2").Value2
$o2 = $ws.Range("O2").Value2
$p2 = $ws.Range("P2").Value2
$q2 = $ws.Range("Q2").Value2
$r2 = $ws.Range("R2").Value2
$s2 = $ws.Range("S2").Value2
$t2 = $ws.Range("T2").Value2

# Save current (pre-edit) values of row 3's changing columns
$d3 = $ws.Range("D3").Value2
$n3 = $ws.Range("N3").Value2
$o3 = $ws.Range("O3").Value2
$p3 = $ws.Range("P3").Value2
$q3 = $ws.Range("Q3").Value2
$r3 = $ws.Range("R3").Value2
$s3 = $ws.Range("S3").Value2
$t3 = $ws.Range("T3").Value2

# Write row 3's former values into row 2
$ws.Range("D2").Value2 = $d3
$ws.Range("N2").Value2 = $n3
$ws.Range("O2").Value2 = $o3
$ws.Range("P2").Value2 = $p3
$ws.Range("Q2").Value2 = $q3
$ws.Range("R2").Value2 = $r3
$ws.Range("S2").Value2 = $s3
$ws.Range("T2").Value2 = $t3

# Write row 2's former values into row 3
$ws.Range("D3").Value2 = $d2
$ws.Range("N3").Value2 = $n2
$ws.Range("O3").Value2 = $o2
$ws.Range("P3").Value2 = $p2
$ws.Range("Q3").Value2 = $q2
$ws.Range("R3").Value2 = $r2
$ws.Range("S3").Value2 = $s2
$ws.Range("T3").Value2 = $t2

$wb.Save()
